$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44607, 900, 1300, 1400, 1350, 1350),
    @(3, 44656, 1000, 900, 1000, 950, 950),
    @(4, 44687, 1000, 1200, 1300, 1250, 1250),
    @(5, 44455, 1100, 900, 1000, 950, 950),
    @(6, 44550, 1300, 1000, 1200, 1100, 1100),
    @(7, 44175, 1600, 1000, 1200, 1100, 1100),
    @(8, 44673, 900, 1300, 1400, 1350, 1350),
    @(9, 44883, 800, 550, 600, 575, 575),
    @(10, 44784, 1000, 1200, 1300, 1250, 1250),
    @(11, 44638, 1000, 900, 950, 925, 925),
    @(12, 44243, 1200, 1200, 1300, 1250, 1250),
    @(13, 44229, 1500, 1400, 1500, 1450, 1450),
    @(14, 44449, 1300, 900, 950, 925, 925),
    @(15, 44341, 1300, 900, 1000, 950, 950),
    @(16, 44407, 1000, 1200, 1300, 1250, 1250),
    @(17, 44291, 1000, 1000, 1200, 1100, 1100),
    @(18, 44649, 600, 900, 1000, 950, 950),
    @(19, 44476, 900, 700, 800, 750, 750),
    @(20, 44453, 1000, 800, 900, 850, 850),
    @(21, 44442, 1250, 850, 900, 875, 875),
    @(22, 44284, 1500, 800, 850, 825, 825),
    @(23, 44484, 900, 750, 800, 775, 775),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P: Precio $/Kg
}
